$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.532.45"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.745.02"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.14"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.744.09"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.57"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.480"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.95"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.365.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.744.86"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.586.09"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.42"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.38"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "500.26"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.15"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.719"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.57"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.83"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.85"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.04%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.03"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000133"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.51%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.46"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.89"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.99"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.35"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.94%  "
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.346"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.05"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +12.12%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "444.02"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.05"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.61%  "
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.44"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.53"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.950.25"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0358"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "139.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.96"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.47"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.10%  "
